$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the gameday rows with literal "Game"/"Practice" values (overwriting the
# formula that previously propagated "Practice" down column B).
$ws.Range("B15").Value = "Game"
$ws.Range("B16").Value = "Practice"
$ws.Range("B21").Value = "Game"
$ws.Range("B22").Value = "Practice"
$ws.Range("B28").Value = "Game"
$ws.Range("B29").Value = "Practice"
$ws.Range("B35").Value = "Game"
$ws.Range("B36").Value = "Practice"
$ws.Range("B42").Value = "Game"
$ws.Range("B43").Value = "Practice"
$ws.Range("B49").Value = "Game"
$ws.Range("B50").Value = "Practice"
$ws.Range("B53").Value = "Game"
$ws.Range("B54").Value = "Practice"
$ws.Range("B56").Value = "Game"
$ws.Range("B57").Value = "Practice"
$ws.Range("B70").Value = "Game"
$ws.Range("B71").Value = "Practice"
$ws.Range("B77").Value = "Game"
$ws.Range("B78").Value = "Practice"
$ws.Range("B84").Value = "Game"
$ws.Range("B85").Value = "Practice"
$ws.Range("B91").Value = "Game"
$ws.Range("B92").Value = "Practice"
$ws.Range("B98").Value = "Game"
$ws.Range("B99").Value = "Practice"

# Restore the view to where the user left off scrolling/selecting.
$ws.Range("B99").Select() | Out-Null
